$wb = $excel.ActiveWorkbook

# The "latest" handback .md URL for ac1a069e-6179-4f35-a566-53815437f1d6 (same
# target already used by the A7 hyperlink on both language sheets).
$latestMdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d0d8a491428e56e61829dc63742041701d9fe9ed/e2e/ac1a069e-6179-4f35-a566-53815437f1d6.md"

# Shared "Error Detail" message (identical text on both language sheets).
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/fcaad4003e49ac746ae24bc350eee15fda75a1e0/e2e/ac1a069e-6179-4f35-a566-53815437f1d6.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d0d8a491428e56e61829dc63742041701d9fe9ed/e2e/ac1a069e-6179-4f35-a566-53815437f1d6.md."

# Excel's "characters" column-width unit is ~5/6 of a character narrower than
# the raw OOXML <col width> value for this font, so asking for 40 - 5/6 here
# lands on an OOXML width of exactly 40.
$colWidthForOoxml40 = 40 - (5 / 6)

function Update-LocalizationSheet {
    param(
        [string]$SheetName,
        [string]$HandbackDateTime,
        [string]$HandbackFile
    )

    $ws = $wb.Worksheets.Item($SheetName)

    # Column P ("Error Detail") widened to fit the new long message.
    $ws.Columns.Item(16).ColumnWidth = $colWidthForOoxml40

    # Row 7 ("ac1a069e-6179-4f35-a566-53815437f1d6") picked up a handback:
    #  I7 - Latest Target File  -> hyperlink to the .md file
    #  J7 - Latest Handback File -> the generated .xlf name
    #  K7 - Latest Handback DateTime -> timestamp of the handback
    #  P7 - Error Detail -> version-mismatch warning
    $ws.Range("J7").Value = $HandbackFile
    $ws.Range("K7").Value = $HandbackDateTime
    $ws.Range("P7").Value = $errorDetail

    $ws.Hyperlinks.Add($ws.Range("I7"), $latestMdUrl, "", "", "ac1a069e-6179-4f35-a566-53815437f1d6.md")
}

Update-LocalizationSheet "zh-cn" "2016-08-23 02:42:43" "ac1a069e-6179-4f35-a566-53815437f1d6.358a283b9cd92c90dbbe51fade2d1d42be3ad461.zh-cn.xlf"
Update-LocalizationSheet "de-de" "2016-08-23 02:42:50" "ac1a069e-6179-4f35-a566-53815437f1d6.358a283b9cd92c90dbbe51fade2d1d42be3ad461.de-de.xlf"
